# Applies the "added more docs to reference" edit to Poster.pptx:
#  - Resizes/repositions the References text box ("Text Box 26")
#  - Replaces its content with two new citation entries, followed by
#    two blank paragraphs (matching the committed OOXML).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)   # "Text Box 26" - the References placeholder

# ---- Resize / reposition the shape -------------------------------------
$shp.Left   = 2304.0000787401573
$shp.Top    = 1554.2274015748033
$shp.Width  = 1125.0462992125983
$shp.Height = 501.01212598425195

# ---- Rebuild the text content -------------------------------------------
$tf = $shp.TextFrame
$tr = $tf.TextRange

$nl = [char]13

$ref1a = 'M. Roland and J. Langer, "Digital Signature Records for the NFC Data Exchange Format," 2010 Second International Workshop on Near Field Communication, Monaco, 2010, pp. 71-76. '
$ref1b = 'doi'
$ref1c = ': 10.1109/NFC.2010.10'
$ref2  = 'electrical 4 u, "Air core transformer," 2011. [Online]. Available: http://www.electrical4u.com/air-core-transformer/. Accessed: Feb. 1, 2017.'

$para1 = $ref1a + $ref1b + $ref1c

$tr.Text = $para1 + $nl + $ref2 + $nl + $nl

# Make sure every paragraph uses the same 24pt size called out in the XML
$tr.Font.Size = 24

# ---- Split paragraph 1 into its three runs -------------------------------
$start1b = $ref1a.Length + 1
$start1c = $ref1a.Length + $ref1b.Length + 1

$run1a = $tr.Characters(1, $ref1a.Length)
$run1a.Font.Size = 24

$run1b = $tr.Characters($start1b, $ref1b.Length)
$run1b.Font.Size = 24

$run1c = $tr.Characters($start1c, $ref1c.Length)
$run1c.Font.Size = 24
